$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# Original layout (consecutive paragraphs):
#   17: "This class represents the expected structure of the response from
#        the Okta API."
#   18: "It has a single property factorResult, which likely contains the
#        result of the OTP verification."
#   19: (already empty paragraph)
#
# Target layout:
#   17: (empty paragraph)
#   19: (already empty paragraph, unchanged)
#
# i.e. paragraph 18 is removed entirely (merging its paragraph mark away)
# and paragraph 17's run content is cleared, leaving an empty paragraph.

$p18 = $null
$p17 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text.TrimEnd("`r")
    if ($t -eq "It has a single property factorResult, which likely contains the result of the OTP verification.") {
        $p18 = $cand
    }
    if ($t -eq "This class represents the expected structure of the response from the Okta API.") {
        $p17 = $cand
    }
}

# Delete paragraph 18 completely (text + its paragraph mark).
$p18.Range.Delete()

# Clear paragraph 17's text but keep the (now empty) paragraph itself --
# delete everything up to (not including) its trailing paragraph mark.
$s = $p17.Range.Start
$e = $p17.Range.End
$d.Range($s, $e - 1).Delete()

# --- Edit 2 ---------------------------------------------------------------
# Paragraph with runs "6." + "Response Handling:" loses the "6." run,
# leaving only "Response Handling:".

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd("`r") -eq "6.Response Handling:") {
        $target = $cand
    }
}

$s2 = $target.Range.Start
$d.Range($s2, $s2 + 2).Delete()
